# POSCON.MECH.BOM.xlsx — Ver 1.2 Rev C -> Ver 1.3 Rev D update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text / label updates -------------------------------------------------
# "Tag:" value (row 1, col F) : POSCON.MECH.v1.2 -> POSCON.MECH.v1.3
$ws.Range("F1").Value = "POSCON.MECH.v1.3"

# Faceplate drawing/part number (row 10, col D)
$ws.Range("D10").Value = "POSCON.FACE VER 1.3 REV D"

# --- "Last Updated:" date (row 3, col F) ----------------------------------
$ws.Range("F3").Value = "9/9/2015"

# --- New blank-but-styled cells in row 17 ---------------------------------
# Materialises A17:C17 as real (empty) cells using the default ("Normal")
# style, mirroring the author manually touching that row in the sheet.
$ws.Range("A17:C17").Borders.LineStyle = -4142

# --- Selection / cursor position ------------------------------------------
$ws.Range("D10").Select()

# --- Cosmetic tab ratio (sheet-tab / scrollbar split) ---------------------
$excel.ActiveWindow.TabRatio = 0.993

# --- Column width tweaks (columns A-H got ~8% wider) -----------------------
# Target (file-stored, "characters") widths are 5.87407407407407,
# 24.0814814814815, 49.5407407407407, 28.6518518518519, 11.8259259259259,
# 19.2444444444444, 14.6740740740741 and 83.5444444444445 for columns A-H.
# This host's ColumnWidth setter snaps to 1/7-character increments
# (MaxDigitWidth=7) before re-adding the fixed 5/7 padding on save, so the
# nearest representable value is used for each column.
$ws.Columns.Item(1).ColumnWidth = 5.142857142857143
$ws.Columns.Item(2).ColumnWidth = 23.428571428571427
$ws.Columns.Item(3).ColumnWidth = 48.857142857142854
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(5).ColumnWidth = 11.142857142857142
$ws.Columns.Item(6).ColumnWidth = 18.571428571428573
$ws.Columns.Item(7).ColumnWidth = 14
$ws.Columns.Item(8).ColumnWidth = 82.85714285714286
